$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPCbS")
$ws.Activate()

# Remove the old explanatory notes in column D (biomass / petroleum notes)
$ws.Range("D9").ClearContents()
$ws.Range("D11").ClearContents()

# Add the three new petroleum-breakdown rows below "petroleum" (row 14)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").NumberFormat = "0"
$ws.Range("B17").Value = 165

# Update the selected cell to reflect the new last used row
[void]$ws.Range("B18").Select()
